# Add a "success" column (D) that flags whether the row's count (column C)
# is within the expected range (0) or is an outlier/success case (1).
# Header D1 gets the same style as the existing headers (B1/C1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-14: default to "0", row 8 (the outlier with count 616) is "1"
# Prefix with an apostrophe so Excel stores the value as text (not a number),
# then reset the style back to Normal so no extra formatting sticks around.
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($r -eq 8) {
        $cell.Value = "'1"
    } else {
        $cell.Value = "'0"
    }
    $cell.Style = "Normal"
}
